$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Populate the two "blank" status rows (10:11 and 12:13) by cloning the
#    formatting of the most recently filled status block (rows 8:9), then
#    filling in the new values for the reporting periods 27.10.2017 and
#    10.11.2017.
# ---------------------------------------------------------------------------
$ws.Range("A8:I9").Copy()
$ws.Range("A10:I13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Merge the cells for both new rows, matching the layout used by every other
# status row in the table.
$ws.Range("A10:B11").Merge()
$ws.Range("C10:E11").Merge()
$ws.Range("F10:F11").Merge()
$ws.Range("G10:G11").Merge()
$ws.Range("H10:H11").Merge()
$ws.Range("I10:I11").Merge()

$ws.Range("A12:B13").Merge()
$ws.Range("C12:E13").Merge()
$ws.Range("F12:F13").Merge()
$ws.Range("G12:G13").Merge()
$ws.Range("H12:H13").Merge()
$ws.Range("I12:I13").Merge()

# Match the row heights of the "thick bottom border" rows (9 -> 11/13).
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(13).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Row 10/11 - week of 27.10.2017 (10% sick)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value2 = 43035

# H10 gets a dedicated text value ("10% krank") instead of a percentage, so
# rebuild its number format / alignment / border to match the rest of the
# row while allowing the cell to hold text.
$h10 = $ws.Range("H10")
$h10.Style = "Normal"
$h10.NumberFormat = "0%"
$h10.HorizontalAlignment = -4152
$h10.Borders.LineStyle = 1
$h10.Borders.Weight = -4138
$h10.Borders.Item(9).LineStyle = -4142
$h10.Font.Name = "Calibri"

$h10.Value2 = "10% krank"
$ws.Range("C10").Value2 = "Versuch Android Studio Problem zu lösen; Research bzgl. Google Maps in C#; Neue Android und C# features"

$ws.Range("F10").Value2 = 0.12
$ws.Range("G10").Value2 = 0.12
$ws.Range("I10").Value2 = 0.12

# ---------------------------------------------------------------------------
# Row 12/13 - week of 10.11.2017
# ---------------------------------------------------------------------------
$ws.Range("A12").Value2 = 43049
$ws.Range("C12").Value2 = "Maps features in C#; Anzeigen von Meetings in Android"

$ws.Range("F12").Value2 = 0.13
$ws.Range("G12").Value2 = 0.13
$ws.Range("H12").Value2 = 0.13
$ws.Range("I12").Value2 = 0.13

# ---------------------------------------------------------------------------
# 2) Misc view/UI bookkeeping that accompanied this status update.
# ---------------------------------------------------------------------------
$ws.Range("C14").Select()
